$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 22:46"

# --- Update Estados Unidos row (row 6): Casos totales, Nuevos casos, Recuperados ---
$ws.Range("B6").Value = 64775
$ws.Range("C6").Value = 9919
$ws.Range("E6").Value = 63472

# --- Update Reino Unido row data (new values) ---
# Reino Unido currently sits in row 13 (country column A13 = "Reino Unido").
# Its updated totals (9529) now exceed Corea del Sur's (9137), so it must move
# above Corea del Sur (to row 12) to keep the sheet sorted descending by
# "Casos totales". Corea del Sur's data itself is unchanged, it simply shifts
# down to row 13.

$ws.Range("A12").Value = "Reino Unido"
$ws.Range("B12").Value = 9529
$ws.Range("C12").Value = 1452
$ws.Range("D12").Value = 135
$ws.Range("E12").Value = 8929
$ws.Range("F12").Value = 163
$ws.Range("G12").Value = 43
$ws.Range("H12").Value = 465

$ws.Range("A13").Value = "Corea del Sur"
$ws.Range("B13").Value = 9137
$ws.Range("C13").Value = 100
$ws.Range("D13").Value = 3730
$ws.Range("E13").Value = 5281
$ws.Range("F13").Value = 59
$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 126
